$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (inline text values).
# A leading apostrophe forces text entry (prevents Excel from auto-
# converting number-looking strings like "238.60" into floating point
# values or applying a text number-format); resetting the style back
# to "Normal" afterwards drops the quote-prefix flag so no stray style
# is left behind on the cell.
$ws.Range("D2").Value = "'35.142.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.91%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.857.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.83%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.20%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'238.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.42%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +1.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'41.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.87%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.330"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.05%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.68%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.05%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.125.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.78%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.93%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'1.834.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.98%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'4.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.91%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'35.114.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.94%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'69.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.45%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'240.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.62%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.49%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +1.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.15%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'168.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.04%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +25.49%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.62%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +2.11%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.43%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.27%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0556"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.56%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +2.49%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +27.23%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.74%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.826"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +18.74%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +12.01%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.48%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +7.47%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'90.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.10%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +4.23%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.340.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.19%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'14.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.02%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.97%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.69%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'12.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +45.29%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +6.46%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.20%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'6.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.98%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.040.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.67%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0679"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.23%  "
$ws.Range("E51").Style = "Normal"
